$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.488.44"
$ws.Range("E2").Value = "  -0.62%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.719.26"
$ws.Range("E3").Value = "  -1.07%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "560.31"
$ws.Range("E5").Value = "  -2.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.88"
$ws.Range("E6").Value = "  +0.74%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.589"
$ws.Range("E8").Value = "  -2.23%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.107"
$ws.Range("E9").Value = "  -2.03%  "
$ws.Range("E10").Value = "  +3.43%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.54"
$ws.Range("E11").Value = "  -1.89%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.371"
$ws.Range("E12").Value = "  -2.59%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.203.57"
$ws.Range("E13").Value = "  -0.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.44"
$ws.Range("E14").Value = "  -0.50%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "63.402.76"
$ws.Range("E15").Value = "  -0.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000147"
$ws.Range("E16").Value = "  -2.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.726.73"
$ws.Range("E17").Value = "  -0.88%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.24"
$ws.Range("E18").Value = "  +0.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.66"
$ws.Range("E19").Value = "  -3.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "351.18"
$ws.Range("E20").Value = "  -1.21%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.49"
$ws.Range("E21").Value = "  -3.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.512"
$ws.Range("E23").Value = "  -3.99%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.73"
$ws.Range("E24").Value = "  -2.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.168"
$ws.Range("E25").Value = "  -1.44%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.19"
$ws.Range("E27").Value = "  -3.91%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0888"
$ws.Range("E28").Value = "  -2.42%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.95"
$ws.Range("E29").Value = "  +0.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.35"
$ws.Range("E30").Value = "  +8.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.10"
$ws.Range("E31").Value = "  -0.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "165.85"
$ws.Range("E32").Value = "  -2.26%  "
$ws.Range("B33").Value = "USDe"
$ws.Range("C33").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.81"
$ws.Range("E34").Value = "  -1.87%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.47"
$ws.Range("E35").Value = "  +0.80%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.82"
$ws.Range("E36").Value = "  -1.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.78"
$ws.Range("E37").Value = "  -0.99%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "339.74"
$ws.Range("E38").Value = "  +1.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.949"
$ws.Range("E39").Value = "  -3.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.19"
$ws.Range("E40").Value = "  +0.27%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.02"
$ws.Range("E41").Value = "  -3.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.38"
$ws.Range("E42").Value = "  -1.64%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.29"
$ws.Range("E43").Value = "  -1.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.56"
$ws.Range("E44").Value = "  -3.91%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0573"
$ws.Range("E45").Value = "  -2.49%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.623"
$ws.Range("E46").Value = "  -0.19%  "
$ws.Range("B47").Value = "FirstDigitalUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.999"
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "131.61"
$ws.Range("E48").Value = "  -2.86%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0246"
$ws.Range("E49").Value = "  -3.38%  "
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0985"
$ws.Range("E50").Value = "  -2.96%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.05"
$ws.Range("E51").Value = "  -0.03%  "
